# Automatic update of files.
#
# 1) Column C ("Förändrad") on every data row moves from serial date
#    45184 (2023-09-15) to 45186 (2023-09-17).
# 2) The HYPERLINK() formulas in columns S, T, V, W, X, Y (where present)
#    gain a second argument — the friendly link text — equal to the
#    row's "Beteckning" value from column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Columns that carry HYPERLINK(...) formulas needing the friendly-text
# second argument.
$linkCols = @(19, 20, 22, 23, 24, 25)   # S, T, V, W, X, Y

for ($r = 2; $r -le $lastRow; $r++) {

    # --- 1) Bump the "Förändrad" date in column C ---------------------
    $cCell = $ws.Cells.Item($r, 3)
    $cVal = $cCell.Value2
    if ($cVal -eq 45184) {
        $cCell.Value = 45186
    }

    # --- 2) Patch HYPERLINK formulas on this row -----------------------
    $designation = $ws.Cells.Item($r, 1).Value2

    foreach ($col in $linkCols) {
        $cell = $ws.Cells.Item($r, $col)
        $formula = $cell.Formula
        if ([string]::IsNullOrEmpty($formula)) {
            continue
        }
        if ($formula -match '^=HYPERLINK\("([^"]*)"\)$') {
            $url = $matches[1]
            $cell.Formula = '=HYPERLINK("' + $url + '", "' + $designation + '")'
        }
    }
}
